$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 408.4
$ws.Range("J19").Value = 385.5
$ws.Range("L19").Value = 385.5
$ws.Range("N19").Value = -735.5
$ws.Range("H97").Value = 220944
$ws.Range("J97").Value = 220944
$ws.Range("L97").Value = 662832
$ws.Range("N97").Value = -663824
$ws.Range("H115").Value = 2605.25
$ws.Range("I115").Value = 2323.2856
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 6969.8568
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -5402.8568
$ws.Range("N115").Value = -12134
$ws.Range("H129").Value = 1273.4
$ws.Range("I129").Value = 2364.4
$ws.Range("J129").Value = 1091.5667
$ws.Range("K129").Value = 7093.200000000001
$ws.Range("L129").Value = 3274.7001
$ws.Range("M129").Value = -2093.200000000001
$ws.Range("N129").Value = -13274.7001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1777
$ws.Range("I61").Value = 1306.6
$ws.Range("J61").Value = 2785
$ws.Range("K61").Value = 1306.6
$ws.Range("L61").Value = 2785
$ws.Range("M61").Value = -1094.6
$ws.Range("N61").Value = -3209
$ws.Range("H136").Value = 1777
$ws.Range("I136").Value = 1306.6
$ws.Range("J136").Value = 2785
$ws.Range("K136").Value = 3919.8
$ws.Range("L136").Value = 8355
$ws.Range("M136").Value = -1369.8
$ws.Range("N136").Value = -13455

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3021.79
$ws.Range("I31").Value = 848.95
$ws.Range("J31").Value = 3565
$ws.Range("K31").Value = 848.95
$ws.Range("L31").Value = 3565
$ws.Range("M31").Value = -553.95
$ws.Range("N31").Value = -4155
$ws.Range("H34").Value = 3021.79
$ws.Range("I34").Value = 848.95
$ws.Range("J34").Value = 3565
$ws.Range("K34").Value = 848.95
$ws.Range("L34").Value = 3565
$ws.Range("M34").Value = -646.95
$ws.Range("N34").Value = -3969
$ws.Range("H50").Value = 39263.332
$ws.Range("J50").Value = 39263.332
$ws.Range("L50").Value = 39263.332
$ws.Range("N50").Value = -40513.332
$ws.Range("H51").Value = 83365976
$ws.Range("J51").Value = 39169.8
$ws.Range("L51").Value = 39169.8
$ws.Range("N51").Value = -40641.8
$ws.Range("H59").Value = 30788.8
$ws.Range("J59").Value = 28460
$ws.Range("L59").Value = 28460
$ws.Range("N59").Value = -30750
$ws.Range("H60").Value = 10847.85
$ws.Range("J60").Value = 11103
$ws.Range("L60").Value = 11103
$ws.Range("N60").Value = -12125
$ws.Range("H61").Value = 83365976
$ws.Range("J61").Value = 39169.8
$ws.Range("L61").Value = 39169.8
$ws.Range("N61").Value = -39865.8
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 397.25
$ws.Range("I12").Value = 335.8
$ws.Range("J12").Value = 413.42105
$ws.Range("K12").Value = 1007.4
$ws.Range("L12").Value = 1240.26315
$ws.Range("M12").Value = -834.4000000000001
$ws.Range("N12").Value = -1586.26315
$ws.Range("H40").Value = 4557.75
$ws.Range("I40").Value = 5688.1113
$ws.Range("J40").Value = 1166.6666
$ws.Range("K40").Value = 22752.4452
$ws.Range("L40").Value = 4666.6664
$ws.Range("M40").Value = -22683.4452
$ws.Range("N40").Value = -4804.6664
$ws.Range("H80").Value = 100200950
$ws.Range("I80").Value = 500625.75
$ws.Range("K80").Value = 1501877.25
$ws.Range("M80").Value = -1500941.25
$ws.Range("H83").Value = 100200950
$ws.Range("I83").Value = 500625.75
$ws.Range("K83").Value = 4505631.75
$ws.Range("M83").Value = -4500951.75
$ws.Range("H102").Value = 16599.666
$ws.Range("I102").Value = 9999
$ws.Range("J102").Value = 19900
$ws.Range("K102").Value = 29997
$ws.Range("L102").Value = 59700
$ws.Range("M102").Value = -27563
$ws.Range("N102").Value = -64568
$ws.Range("H114").Value = 3171.6
$ws.Range("I114").Value = 1056.75
$ws.Range("J114").Value = 4581.5
$ws.Range("K114").Value = 3170.25
$ws.Range("L114").Value = 13744.5
$ws.Range("M114").Value = 83.75
$ws.Range("N114").Value = -20252.5
$ws.Range("H123").Value = 2000
$ws.Range("I123").Value = 2000
$ws.Range("K123").Value = 6000
$ws.Range("M123").Value = -3550
$ws.Range("H131").Value = 3553.1396
$ws.Range("I131").Value = 12876.125
$ws.Range("J131").Value = 1422.1714
$ws.Range("K131").Value = 38628.375
$ws.Range("L131").Value = 4266.5142
$ws.Range("M131").Value = -33588.375
$ws.Range("N131").Value = -14346.5142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 46243.5
$ws.Range("J110").Value = 46243.5
$ws.Range("L110").Value = 46243.5
$ws.Range("N110").Value = -54423.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3038.625
$ws.Range("I7").Value = 2050.6667
$ws.Range("J7").Value = 6002.5
$ws.Range("K7").Value = 2050.6667
$ws.Range("L7").Value = 6002.5
$ws.Range("M7").Value = -1938.6667
$ws.Range("N7").Value = -6226.5
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 666.6667
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 666.6667
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = -1256.6667
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 666.6667
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 666.6667
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = -880.6667
$ws.Range("H111").Value = 46116.75
$ws.Range("J111").Value = 46116.75
$ws.Range("L111").Value = 46116.75
$ws.Range("N111").Value = -54296.75
$ws.Range("H126").Value = 3038.625
$ws.Range("I126").Value = 2050.6667
$ws.Range("J126").Value = 6002.5
$ws.Range("K126").Value = 6152.000100000001
$ws.Range("L126").Value = 18007.5
$ws.Range("M126").Value = -3682.000100000001
$ws.Range("N126").Value = -22947.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1996.6
$ws.Range("I126").Value = 1996.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5989.799999999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3519.799999999999
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 556662.4399999999
$ws.Range("I136").Value = 667521.8
$ws.Range("K136").Value = 2002565.4
$ws.Range("M136").Value = -2000015.4
